$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 holds account 005002457 / ROSANGELA / 17000 - remove it entirely,
# shifting all subsequent rows up.
$ws.Rows.Item(4).Delete()

# The RAPHAELA row (005366255) is now row 4; update its Saldo from 1700 to 1000.
$ws.Cells.Item(4, 3).Value = 1000
